# "Add files via upload" - re-upload of the workbook after filling in the
# missing "Survey 3" row (row 4) data on Sheet1, and updating the active
# selection left on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 25
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 0

$ws.Range("E14").Select()
